# Dodanie podziału treningu na części
#
# - Adds a new "Trening" column (F) classifying each measurement row into
#   "Duża Gra" or "Mała Gra".
# - Converts the Timestamp column (A) from text strings into real Excel
#   date-time serial values formatted with a custom "YYYY-MM-DD HH:MM:SS"
#   number format.
# - Corrects row 5's Seconds/Velocity/Acceleration_SMA values.
# - Appends 6 additional measurement rows (8-13) belonging to "Mała Gra".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header for the new column F (reuse the bold/bordered header style) ---
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# --- Row data ---
# Columns: TimestampSerial, Seconds, Velocity, Acceleration_SMA, Velocity_Bin, Trening
$rows = @(
    @{R=2;  A=45685.64817673611; B=1117.4; C=13.82;              D=3.76764794758388;  E="10-15"; F="Duża Gra"},
    @{R=3;  A=45685.65527627315; B=1730.8; C=13.3;                D=3.434327329908097; E="10-15"; F="Duża Gra"},
    @{R=4;  A=45685.65701122685; B=1880.7; C=14.99;               D=3.497796637671334; E="10-15"; F="Duża Gra"},
    @{R=5;  A=45685.64689085648; B=1006.3; C=9.130000000000001;   D=2.99323902811323;  E="5-10";  F="Duża Gra"},
    @{R=6;  A=45685.64817326389; B=1117.1; C=9.08;                D=3.41945140702384;  E="5-10";  F="Duża Gra"},
    @{R=7;  A=45685.64952974537; B=1234.3; C=9.26;                D=3.240590572357179; E="5-10";  F="Duża Gra"},
    @{R=8;  A=45685.66956909722; B=2965.7; C=14.57;               D=3.408025537218367; E="10-15"; F="Mała Gra"},
    @{R=9;  A=45685.67857372685; B=3743.7; C=14.93;               D=3.41138824394771;  E="10-15"; F="Mała Gra"},
    @{R=10; A=45685.6805181713;  B=3911.7; C=13.87;               D=3.357474974223547; E="10-15"; F="Mała Gra"},
    @{R=11; A=45685.66956446759; B=2965.3; C=9.130000000000001;   D=2.642762865339007; E="5-10";  F="Mała Gra"},
    @{R=12; A=45685.67559456018; B=3486.3; C=9.52;                D=2.744750993592396; E="5-10";  F="Mała Gra"},
    @{R=13; A=45685.67892789352; B=3774.3; C=9.94;                D=2.544219238417492; E="5-10";  F="Mała Gra"}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}

# --- Apply the date/time number format to the Timestamp column ---
# First apply the lowercase variant to a single cell (registers numFmt 164),
# then apply the actual uppercase variant to the whole range (registers and
# uses numFmt 165), matching how the format ended up recorded in the file.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
